# The workbook groups cleaned article filenames by topic column
# (A=Tech, B=Sports/overflow, C=Business-Politics, D=Food, E=Science).
# A handful of articles had been misfiled under the wrong topic; this
# reclassifies them into the right column. Because each column is kept
# sorted/compacted with no gaps, moving an entry also re-packs every
# following cell in both the source and destination columns by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = 'cleaned_0308_tech'
$ws.Range("B28").Value = 'cleaned_0309_tech'
$ws.Range("B29").Value = 'cleaned_0318_tech'
$ws.Range("B30").Value = 'cleaned_0504_business'
$ws.Range("B31").Value = 'cleaned_0505_politics'
$ws.Range("B32").Value = 'cleaned_0506_business'
$ws.Range("B34").Value = 'cleaned_0511_business'
$ws.Range("B35").Value = 'cleaned_0513_business'
$ws.Range("B36").Value = 'cleaned_0515_politics'
$ws.Range("B37").Value = 'cleaned_0516_business'
$ws.Range("B38").Value = 'cleaned_0518_business'
$ws.Range("B39").Value = 'cleaned_0519_business'
$ws.Range("B40").Value = 'cleaned_0520_politics'
$ws.Range("B41").Value = 'cleaned_9901_sports'
$ws.Range("B42").Value = 'cleaned_9903_tech'
$ws.Range("B43").Value = 'cleaned_9905_business'
$ws.Range("B44").Value = 'cleaned_9906_politics'
$ws.Range("B45").Value = 'cleaned_science_0409'
$ws.Range("C15").Value = 'cleaned_0508_business'
$ws.Range("C16").Value = 'cleaned_0509_business'
$ws.Range("C17").Value = 'cleaned_0514_business'
$ws.Range("C18").Value = 'cleaned_0517_business'
$ws.Range("C19").Value = 'cleaned_tech_0401'
$ws.Range("C20").Value = 'cleaned_tech_0402'
$ws.Range("C21").Value = 'cleaned_tech_0403'
$ws.Range("C22").Value = 'cleaned_tech_0404'
$ws.Range("D6").Value = 'cleaned_0202_food'
$ws.Range("D7").Value = 'cleaned_0204_food'
$ws.Range("D8").Value = 'cleaned_0207_food'
$ws.Range("D9").Value = 'cleaned_0208_food'
$ws.Range("D10").Value = 'cleaned_0211_food'
$ws.Range("D11").Value = 'cleaned_0217_food'
$ws.Range("D12").Value = 'cleaned_0218_food'
$ws.Range("D13").Value = 'cleaned_0219_food'
$ws.Range("D14").Value = 'cleaned_0306_tech'
$ws.Range("D15").Value = 'cleaned_0507_business'
$ws.Range("D16").Value = 'cleaned_9902_food'
$ws.Range("E2").Value = 'cleaned_0206_food'
$ws.Range("E3").Value = 'cleaned_0305_science'
$ws.Range("E6").Value = 'cleaned_9904_science'
$ws.Range("E7").Value = 'cleaned_science_0401'
$ws.Range("E8").Value = 'cleaned_science_0402'
$ws.Range("E9").Value = 'cleaned_science_0403'
$ws.Range("E10").Value = 'cleaned_science_0404'
$ws.Range("E11").Value = 'cleaned_science_0405'
$ws.Range("E12").Value = 'cleaned_science_0406'
$ws.Range("E13").Value = 'cleaned_science_0407'
$ws.Range("E14").Value = 'cleaned_science_0408'
$ws.Range("E15").Value = 'cleaned_science_0410'
$ws.Range("E16").Value = 'cleaned_science_0411'
$ws.Range("E17").Value = 'cleaned_science_0412'
$ws.Range("E18").Value = 'cleaned_science_0413'
$ws.Range("E19").Value = 'cleaned_science_0414'
$ws.Range("E20").Value = 'cleaned_science_0415'

$ws.Range("C23").ClearContents()
$ws.Range("C24").ClearContents()
$ws.Range("C25").ClearContents()
$ws.Range("C26").ClearContents()
